# Auto-generated Excel COM-interop edit script
# Updates cryptocurrency price/volume table rows 2-51 per the commit's data refresh,
# and fixes two row-ordering swaps (Stellar/VeChain, Quant/Cronos).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "23.328.92"
$ws.Range("E2").Value = "  -1.64%  "

# Row 3
$ws.Range("D3").Value = "1.629.02"
$ws.Range("E3").Value = "  -1.77%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  +0.00%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "297.95"
$ws.Range("E6").Value = "  -1.66%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3765"
$ws.Range("E7").Value = "  -1.36%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.14"
$ws.Range("E8").Value = "  -2.18%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3463"
$ws.Range("E9").Value = "  -4.16%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08029"
$ws.Range("E10").Value = "  -1.99%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.200"
$ws.Range("E11").Value = "  -2.40%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.04%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.87"
$ws.Range("E13").Value = "  -2.93%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.289"
$ws.Range("E14").Value = "  -2.81%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.218"
$ws.Range("E15").Value = "  -2.55%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001187"
$ws.Range("E16").Value = "  -3.28%  "

# Row 17
$ws.Range("D17").Value = "1.623.02"
$ws.Range("E17").Value = "  -1.64%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.55"
$ws.Range("E18").Value = "  -3.45%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06939"
$ws.Range("E19").Value = "  -1.06%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.598"
$ws.Range("E20").Value = "  -3.28%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.24"
$ws.Range("E21").Value = "  -2.20%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.08%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.33"
$ws.Range("E23").Value = "  -3.91%  "

# Row 24
$ws.Range("D24").Value = "23.329.08"
$ws.Range("E24").Value = "  -1.61%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.428"
$ws.Range("E25").Value = "  -3.08%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.955"
$ws.Range("E26").Value = "  -1.20%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.89"
$ws.Range("E27").Value = "  -1.57%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.62"
$ws.Range("E28").Value = "  -1.35%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.159"
$ws.Range("E29").Value = "  -1.15%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.57"
$ws.Range("E30").Value = "  -2.07%  "

# Row 31
$ws.Range("D31").Value = "1.805.47"
$ws.Range("E31").Value = "  -1.71%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.701"
$ws.Range("E32").Value = "  -5.50%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.116"
$ws.Range("E33").Value = "  -4.99%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.16"
$ws.Range("E34").Value = "  -7.07%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9748"
$ws.Range("E35").Value = "  -7.65%  "

# Row 36
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02654"
$ws.Range("E36").Value = "  -5.64%  "

# Row 37
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08740"
$ws.Range("E37").Value = "  -0.68%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2415"
$ws.Range("E38").Value = "  -3.94%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.833"
$ws.Range("E39").Value = "  -4.12%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06682"
$ws.Range("E40").Value = "  -4.63%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.68"
$ws.Range("E41").Value = "  -2.36%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6781"
$ws.Range("E42").Value = "  -3.13%  "

# Row 43
$ws.Range("E43").Value = "  -3.04%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.33"
$ws.Range("E44").Value = "  -4.28%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.03%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6298"
$ws.Range("E46").Value = "  -3.32%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.228"
$ws.Range("E47").Value = "  -3.55%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.879"
$ws.Range("E48").Value = "  -2.16%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07623"
$ws.Range("E49").Value = "  -3.73%  "

# Row 50
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "126.05"
$ws.Range("E50").Value = "  -1.69%  "

# Row 51
$ws.Range("E51").Value = "  +1.43%  "
